$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows with new daily power record data
$ws.Range("C63").Value = 0.24027777777777778

$ws.Range("B64").Value = 0
$ws.Range("C64").Value = 0

$ws.Range("B65").Value = 0.73958333333333337
$ws.Range("C65").Value = 0.99930555555555556

# Add new row 66 for 2018-10-17
$ws.Range("A66").Value = 43390
$ws.Range("B66").Value = 0

$ws.Range("D66").Formula = "=(C66-B66)* 1440"
$ws.Range("E66").Formula = "=IF(C66>B66, (C66-B66)*1440, (B66-C66)*1440)"
$ws.Range("F66").Formula = "=ABS((C66-B66)*1440)"

# Match the number formatting used by the row above for the new row
$ws.Range("D65:F65").Copy()
$ws.Range("D66:F66").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Grow the table to include the new row
$tbl = $ws.ListObjects.Item("comforter_cda_table")
$tbl.Resize($ws.Range("A1:F66"))

# Update selection/view state to reflect where editing left off
$ws.Range("C66").Select()
